$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $orig = $range.Style
    $range.Value = "'" + $value
    $range.Style = $orig
}

Set-TextValue $ws.Range("D2") "63.665.39"
Set-TextValue $ws.Range("E2") "  -1.34%  "
Set-TextValue $ws.Range("D3") "3.041.79"
Set-TextValue $ws.Range("E3") "  -1.58%  "
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "555.31"
Set-TextValue $ws.Range("E5") "  -0.41%  "
Set-TextValue $ws.Range("D6") "141.27"
Set-TextValue $ws.Range("E6") "  -1.87%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("D8") "3.039.48"
Set-TextValue $ws.Range("E8") "  -1.47%  "
Set-TextValue $ws.Range("E9") "  +3.81%  "
Set-TextValue $ws.Range("E10") "  +0.15%  "
Set-TextValue $ws.Range("D11") "6.16"
Set-TextValue $ws.Range("E11") "  -13.54%  "
Set-TextValue $ws.Range("E12") "  +4.37%  "
Set-TextValue $ws.Range("E13") "  +0.09%  "
Set-TextValue $ws.Range("D14") "35.39"
Set-TextValue $ws.Range("E14") "  +0.29%  "
Set-TextValue $ws.Range("D15") "3.536.72"
Set-TextValue $ws.Range("E15") "  -1.67%  "
Set-TextValue $ws.Range("D16") "63.677.30"
Set-TextValue $ws.Range("E16") "  -1.40%  "
Set-TextValue $ws.Range("D17") "3.033.49"
Set-TextValue $ws.Range("E17") "  -1.84%  "
Set-TextValue $ws.Range("E18") "  +0.40%  "
Set-TextValue $ws.Range("D19") "6.74"
Set-TextValue $ws.Range("E19") "  -0.76%  "
Set-TextValue $ws.Range("D20") "472.52"
Set-TextValue $ws.Range("E20") "  -1.96%  "
Set-TextValue $ws.Range("D21") "14.02"
Set-TextValue $ws.Range("E21") "  +1.74%  "
Set-TextValue $ws.Range("D22") "0.680"
Set-TextValue $ws.Range("E22") "  +0.67%  "
Set-TextValue $ws.Range("D23") "14.45"
Set-TextValue $ws.Range("E23") "  +8.41%  "
Set-TextValue $ws.Range("D24") "7.51"
Set-TextValue $ws.Range("E24") "  -0.51%  "
Set-TextValue $ws.Range("E25") "  +1.89%  "
Set-TextValue $ws.Range("E26") "  +0.06%  "
Set-TextValue $ws.Range("D27") "2.78"
Set-TextValue $ws.Range("E27") "  +0.32%  "
Set-TextValue $ws.Range("D28") "8.05"
Set-TextValue $ws.Range("E28") "  -1.47%  "
Set-TextValue $ws.Range("D29") "2.02"
Set-TextValue $ws.Range("E29") "  -1.98%  "
Set-TextValue $ws.Range("E30") "  -0.13%  "
Set-TextValue $ws.Range("D31") "26.03"
Set-TextValue $ws.Range("E31") "  +0.14%  "
Set-TextValue $ws.Range("E32") "  -0.32%  "
Set-TextValue $ws.Range("D33") "2.43"
Set-TextValue $ws.Range("E33") "  -0.81%  "
Set-TextValue $ws.Range("D34") "5.66"
Set-TextValue $ws.Range("E34") "  -0.26%  "
Set-TextValue $ws.Range("D35") "6.18"
Set-TextValue $ws.Range("E35") "  -0.46%  "
Set-TextValue $ws.Range("D36") "54.88"
Set-TextValue $ws.Range("E36") "  +0.17%  "
Set-TextValue $ws.Range("D37") "0.0407"
Set-TextValue $ws.Range("E37") "  +0.09%  "
Set-TextValue $ws.Range("D38") "438.06"
Set-TextValue $ws.Range("E38") "  -5.76%  "
Set-TextValue $ws.Range("D39") "0.0810"
Set-TextValue $ws.Range("E39") "  -1.69%  "
Set-TextValue $ws.Range("D40") "2.995.78"
Set-TextValue $ws.Range("E40") "  -0.26%  "
Set-TextValue $ws.Range("D41") "2.73"
Set-TextValue $ws.Range("E41") "  -4.04%  "
Set-TextValue $ws.Range("B42") "Cosmos"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D42") "8.25"
Set-TextValue $ws.Range("E42") "  -0.03%  "
Set-TextValue $ws.Range("B43") "Kaspa"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D43") "0.116"
Set-TextValue $ws.Range("E43") "  +0.23%  "
Set-TextValue $ws.Range("D44") "0.268"
Set-TextValue $ws.Range("E44") "  +3.96%  "
Set-TextValue $ws.Range("D45") "27.69"
Set-TextValue $ws.Range("E45") "  -1.68%  "
Set-TextValue $ws.Range("D46") "2.23"
Set-TextValue $ws.Range("E46") "  +6.13%  "
Set-TextValue $ws.Range("E47") "  -0.05%  "
Set-TextValue $ws.Range("E48") "  +0.09%  "
Set-TextValue $ws.Range("D49") "117.77"
Set-TextValue $ws.Range("E49") "  +0.03%  "
Set-TextValue $ws.Range("D50") "0.0₃0509"
Set-TextValue $ws.Range("E50") "  -1.18%  "
Set-TextValue $ws.Range("E51") "  +0.12%  "
